# Updates cached FFXIV market-board prices and recomputed Leve profit figures
# for the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR crafting tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4327.1113
$ws.Range("I70").Value = 1245
$ws.Range("J70").Value = 4712.375
$ws.Range("K70").Value = 3735
$ws.Range("L70").Value = 14137.125
$ws.Range("M70").Value = -3465
$ws.Range("N70").Value = -14677.125
$ws.Range("H73").Value = 4327.1113
$ws.Range("I73").Value = 1245
$ws.Range("J73").Value = 4712.375
$ws.Range("K73").Value = 3735
$ws.Range("L73").Value = 14137.125
$ws.Range("M73").Value = -2799
$ws.Range("N73").Value = -16009.125
$ws.Range("H86").Value = 1399
$ws.Range("I86").Value = 1399
$ws.Range("K86").Value = 1399
$ws.Range("M86").Value = -276
$ws.Range("H87").Value = 58284
$ws.Range("J87").Value = 69855
$ws.Range("L87").Value = 69855
$ws.Range("N87").Value = -72351
$ws.Range("H89").Value = 1399
$ws.Range("I89").Value = 1399
$ws.Range("K89").Value = 6995
$ws.Range("M89").Value = -1379
$ws.Range("H90").Value = 58284
$ws.Range("J90").Value = 69855
$ws.Range("L90").Value = 209565
$ws.Range("N90").Value = -222045
$ws.Range("H112").Value = 3161.7144
$ws.Range("J112").Value = 3161.7144
$ws.Range("L112").Value = 9485.143199999999
$ws.Range("N112").Value = -11701.1432
$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 4000
$ws.Range("K125").Value = 36000
$ws.Range("M125").Value = -33540
$ws.Range("H127").Value = 4473.5
$ws.Range("I127").Value = 4523
$ws.Range("J127").Value = 4325
$ws.Range("K127").Value = 13569
$ws.Range("L127").Value = 12975
$ws.Range("M127").Value = -8609
$ws.Range("N127").Value = -22895
$ws.Range("H141").Value = 2370.1
$ws.Range("I141").Value = 1300.5714
$ws.Range("J141").Value = 4865.6665
$ws.Range("K141").Value = 3901.7142
$ws.Range("L141").Value = 14596.9995
$ws.Range("M141").Value = 1278.2858
$ws.Range("N141").Value = -24956.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 99.76667
$ws.Range("I5").Value = 102.333336
$ws.Range("J5").Value = 99.125
$ws.Range("K5").Value = 102.333336
$ws.Range("L5").Value = 99.125
$ws.Range("M5").Value = 9.666663999999997
$ws.Range("N5").Value = -323.125
$ws.Range("H32").Value = 7569.6665
$ws.Range("I32").Value = 7569.6665
$ws.Range("K32").Value = 7569.6665
$ws.Range("M32").Value = -7282.6665
$ws.Range("H110").Value = 2820.4614
$ws.Range("I110").Value = 1161.5385
$ws.Range("K110").Value = 1161.5385
$ws.Range("M110").Value = 883.4614999999999
$ws.Range("H132").Value = 3721.8572
$ws.Range("I132").Value = 3113.6316
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 9340.8948
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -6810.8948
$ws.Range("N132").Value = -33560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 99.76667
$ws.Range("I4").Value = 102.333336
$ws.Range("J4").Value = 99.125
$ws.Range("K4").Value = 102.333336
$ws.Range("L4").Value = 99.125
$ws.Range("M4").Value = 12.666664
$ws.Range("N4").Value = -329.125
$ws.Range("H81").Value = 97831.5
$ws.Range("J81").Value = 97831.5
$ws.Range("L81").Value = 97831.5
$ws.Range("N81").Value = -99953.5
$ws.Range("H84").Value = 97831.5
$ws.Range("J84").Value = 97831.5
$ws.Range("L84").Value = 293494.5
$ws.Range("N84").Value = -304102.5
$ws.Range("H105").Value = 1662.6666
$ws.Range("J105").Value = 1398
$ws.Range("L105").Value = 1398
$ws.Range("N105").Value = -4892

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.30768999999999
$ws.Range("I7").Value = 90.75
$ws.Range("K7").Value = 90.75
$ws.Range("M7").Value = 22.25
$ws.Range("H12").Value = 649.44446
$ws.Range("I12").Value = 241.66667
$ws.Range("J12").Value = 1465
$ws.Range("K12").Value = 241.66667
$ws.Range("L12").Value = 1465
$ws.Range("M12").Value = -71.66667000000001
$ws.Range("N12").Value = -1805
$ws.Range("H31").Value = 7142.933
$ws.Range("I31").Value = 3183.0908
$ws.Range("J31").Value = 9435.474
$ws.Range("K31").Value = 3183.0908
$ws.Range("L31").Value = 9435.474
$ws.Range("M31").Value = -2888.0908
$ws.Range("N31").Value = -10025.474
$ws.Range("H34").Value = 7142.933
$ws.Range("I34").Value = 3183.0908
$ws.Range("J34").Value = 9435.474
$ws.Range("K34").Value = 3183.0908
$ws.Range("L34").Value = 9435.474
$ws.Range("M34").Value = -2981.0908
$ws.Range("N34").Value = -9839.474
$ws.Range("H107").Value = 528.5333000000001
$ws.Range("I107").Value = 520.2857
$ws.Range("J107").Value = 644
$ws.Range("K107").Value = 520.2857
$ws.Range("L107").Value = 644
$ws.Range("M107").Value = 1399.7143
$ws.Range("N107").Value = -4484
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2004
$ws.Range("I51").Value = 2004
$ws.Range("K51").Value = 6012
$ws.Range("M51").Value = -5552
$ws.Range("H104").Value = 7266.6665
$ws.Range("H112").Value = 2027
$ws.Range("I112").Value = 2027
$ws.Range("K112").Value = 6081
$ws.Range("M112").Value = -4973
$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8991
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6541
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 11287.333
$ws.Range("I139").Value = 11914.5
$ws.Range("K139").Value = 35743.5
$ws.Range("M139").Value = -30603.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("N49").Value = -18368
$ws.Range("H80").Value = 3999.5
$ws.Range("I80").Value = 3999.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3999.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3001.5
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3999.5
$ws.Range("I83").Value = 3999.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 19997.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -15005.5
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 3551.05
$ws.Range("I102").Value = 2991.4
$ws.Range("K102").Value = 2991.4
$ws.Range("M102").Value = -1369.4
$ws.Range("H126").Value = 6013.8
$ws.Range("I126").Value = 6013
$ws.Range("K126").Value = 18039
$ws.Range("M126").Value = -15569

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 208799.8
$ws.Range("I20").Value = 9000
$ws.Range("K20").Value = 9000
$ws.Range("M20").Value = -8774
$ws.Range("H40").Value = 6275.273
$ws.Range("I40").Value = 6096
$ws.Range("K40").Value = 6096
$ws.Range("M40").Value = -5960
$ws.Range("H42").Value = 14497.5
$ws.Range("J42").Value = 19000
$ws.Range("L42").Value = 19000
$ws.Range("N42").Value = -20126
$ws.Range("H46").Value = 7607.9287
$ws.Range("J46").Value = 5896.75
$ws.Range("L46").Value = 5896.75
$ws.Range("N46").Value = -6272.75
$ws.Range("H49").Value = 14497.5
$ws.Range("J49").Value = 19000
$ws.Range("L49").Value = 19000
$ws.Range("N49").Value = -19294
$ws.Range("H61").Value = 3622.6843
$ws.Range("I61").Value = 2455.4
$ws.Range("K61").Value = 2455.4
$ws.Range("M61").Value = -2253.4
$ws.Range("H82").Value = 6090
$ws.Range("H85").Value = 6090
$ws.Range("H92").Value = 100000
$ws.Range("I92").Value = 100000
$ws.Range("K92").Value = 100000
$ws.Range("M92").Value = -97504
$ws.Range("H113").Value = 3622.6843
$ws.Range("I113").Value = 2455.4
$ws.Range("K113").Value = 2455.4
$ws.Range("M113").Value = -285.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8268
$ws.Range("J62").Value = 8457
$ws.Range("L62").Value = 8457
$ws.Range("N62").Value = -9705
$ws.Range("H65").Value = 8268
$ws.Range("J65").Value = 8457
$ws.Range("L65").Value = 42285
$ws.Range("N65").Value = -48525
$ws.Range("H107").Value = 670.3333
$ws.Range("I107").Value = 511.14285
$ws.Range("J107").Value = 1227.5
$ws.Range("K107").Value = 1533.42855
$ws.Range("L107").Value = 3682.5
$ws.Range("M107").Value = 386.5714499999999
$ws.Range("N107").Value = -7522.5
$ws.Range("H126").Value = 2616.7144
$ws.Range("I126").Value = 1551.95
$ws.Range("J126").Value = 5278.625
$ws.Range("K126").Value = 4655.85
$ws.Range("L126").Value = 15835.875
$ws.Range("M126").Value = -2185.85
$ws.Range("N126").Value = -20775.875
$ws.Range("H132").Value = 3802.92
$ws.Range("I132").Value = 3671.2
$ws.Range("J132").Value = 4329.8
$ws.Range("K132").Value = 11013.6
$ws.Range("L132").Value = 12989.4
$ws.Range("M132").Value = -8483.599999999999
$ws.Range("N132").Value = -18049.4
$ws.Range("H135").Value = 40212.875
$ws.Range("J135").Value = 40212.875
$ws.Range("L135").Value = 40212.875
$ws.Range("N135").Value = -50352.875
$ws.Range("H136").Value = 3230.9355
$ws.Range("I136").Value = 2007.75
$ws.Range("K136").Value = 6023.25
$ws.Range("M136").Value = -3473.25

Write-Host "Applied all cell updates."
